# Update countries & provincias Spain
#
# Refreshes the Covid "Pais" sheet (Casos totales / Nuevos casos / Casos
# activos / Recuperados / Casos criticos / Muertes hoy / Muertes, columns
# B-H) with newer counters, and re-applies the descending-by-total-cases
# sort order. Most rows just get new numbers, but a handful of countries
# that were neck-and-neck on total cases swap ranks with their neighbour,
# so those rows also get the country name (column A) swapped along with
# the new figures:
#   - row 82/83   Bosnia y Herzegovina  <-> Serbia
#   - row 121/122 Guadalupe             <-> Sri Lanka
#   - row 216/217 Islas Malvinas        <-> Montserrat

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Datos actualizados a ..." banner
$ws.Range("A1").Value = "Datos actualizados a 24 de Octubre de 2020 a las 16:19"

$cols = @("B", "C", "D", "E", "F", "G", "H")

function Set-RowData($row, $values) {
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$row").Value = $values[$i]
    }
}

# Estados Unidos
Set-RowData 4 @(8756462, 9509, 5699030, 2828056, 0, 92, 229376)

# India
Set-RowData 5 @(7815420, 1752, 7016046, 681365, 0, 17, 118009)

# Chile
Set-RowData 17 @(500542, 1631, 476576, 10074, 0, 48, 13892)

# Irak
Set-RowData 19 @(449153, 3204, 378209, 60376, 0, 55, 10568)

# Alemania
Set-RowData 20 @(421101, 3751, 314100, 96908, 0, 3, 10093)

# Portugal
Set-RowData 44 @(116109, 3669, 67842, 45970, 0, 21, 2297)

# Bielorrusia
$ws.Range("G54").Value = 4

# Estado de Palestina
Set-RowData 72 @(49989, 410, 43232, 6314, 0, 4, 443)

# Ghana
Set-RowData 77 @(47690, 89, 46887, 487, 0, 2, 316)

# Birmania
Set-RowData 78 @(43788, 1423, 23708, 19014, 0, 28, 1066)

# Row 82 becomes Serbia (was Bosnia y Herzegovina)
$ws.Range("A82").Value = "Serbia"
Set-RowData 82 @(38872, 757, 31536, 6547, 0, 3, 789)

# Row 83 becomes Bosnia y Herzegovina (was Serbia)
$ws.Range("A83").Value = "Bosnia y Herzegovina"
Set-RowData 83 @(38493, 0, 26260, 11168, 0, 0, 1065)

# Noruega
Set-RowData 97 @(17670, 138, 11863, 5528, 0, 0, 279)

# Uganda
Set-RowData 109 @(11297, 134, 7281, 3917, 0, 0, 99)

# Row 121 becomes Sri Lanka (was Guadalupe)
$ws.Range("A121").Value = "Sri Lanka"
Set-RowData 121 @(7354, 201, 3714, 3625, 0, 1, 15)

# Row 122 becomes Guadalupe (was Sri Lanka)
$ws.Range("A122").Value = "Guadalupe"
Set-RowData 122 @(7329, 0, 2199, 5015, 0, 0, 115)

# Lesoto
Set-RowData 163 @(1940, 6, 970, 927, 0, 0, 43)

# Row 216 becomes Montserrat (was Islas Malvinas)
$ws.Range("A216").Value = "Montserrat"
Set-RowData 216 @(13, 0, 12, 0, 0, 0, 1)

# Row 217 becomes Islas Malvinas (was Montserrat)
$ws.Range("A217").Value = "Islas Malvinas"
Set-RowData 217 @(13, 0, 13, 0, 0, 0, 0)
